$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "Notice u/s 94 BNSS, 2023" paragraph (the large,
#    centered, italic heading that currently sits right above the
#    "Subject: -" paragraph).
# ------------------------------------------------------------------
$findRange = $d.Content
$oldFound = $findRange.Find.Execute("Notice u/s 94 BNSS, 2023", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($oldFound) {
    $oldParaIndex = $findRange.Paragraphs.Item(1).Index
    $oldParaRange = $d.Paragraphs.Item($oldParaIndex).Range
    $oldParaRange.Delete()
}

# ------------------------------------------------------------------
# 2) Insert a new "Notice u/s 94 BNSS, 2023" paragraph (bold,
#    underlined, justified) immediately above the "To," paragraph.
# ------------------------------------------------------------------
$toRange = $d.Content
$toFound = $toRange.Find.Execute("To,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($toFound) {
    $toParaIndex = $toRange.Paragraphs.Item(1).Index
    $toParaRange = $d.Paragraphs.Item($toParaIndex).Range
    $newPara = $toParaRange.InsertParagraphBefore()

    $newParaIndex = $toParaIndex
    $newRange = $d.Paragraphs.Item($newParaIndex).Range

    $noticeXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Notice u/s 94 BNSS, 2023</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $newRange.InsertXML($noticeXml)
}
